# Workbook already open as $excel.ActiveWorkbook
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("UserList")
$ws2 = $wb.Worksheets.Item("Sorting")

# Find/replace typo'd product names in the Sorting sheet's shared strings
# ("...OnesieM" / "...BackpackM" -> the correct names without the trailing "M")
$ws2.Range("B3").Value = "Sauce Labs Onesie"
$ws2.Range("B5").Value = "Sauce Labs Backpack"

# Duplicate the UserList sheet (Excel names the copy "UserList (2)" and
# places it after the last existing sheet) for a smoke-test data set
$ws1.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count)) | Out-Null
$ws3 = $wb.Worksheets.Item("UserList (2)")

# Remove the locked_out_user row from the new copy - not needed for the smoke test
$ws3.Rows.Item(3).Delete() | Out-Null
$ws3.Rows.Item(3).Select() | Out-Null

# Leave the cursor parked where the author left off on each sheet
$ws1.Range("A17").Select() | Out-Null
$ws2.Range("D15").Select() | Out-Null
